$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, copying the style from the existing
# last header cell (AC1) so the new cells match the bold/bordered header
# formatting used throughout row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins / Losses / Ties) for every player row.
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
